$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest figures.
# D-column text values that look numeric must be forced to remain text (matching
# the original inlineStr cell type) without leaving a residual cell style, so we
# briefly mark the cell as Text, assign the value, then restore the default style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.309.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.077.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.30%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0832"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.39%  "

$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.385.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.778"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.077.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.268.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0828"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.70%  "

$ws.Range("E30").Value = "  +7.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.92%  "

$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("E36").Value = "  +1.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.87%  "

$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.538.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.16%  "

$ws.Range("E43").Value = "  +2.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0919"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.32%  "

$ws.Range("E47").Value = "  +0.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "

$ws.Range("E49").Value = "  +2.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.278.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
